$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: fill in previously-empty Accuracy/Specificity/Sensitivity values ---
$ws.Range("G8").Value = 0.9652
$ws.Range("H8").Value = 0.9722
$ws.Range("I8").Value = 0.9628

# --- Row 9: new "spacer" row (A:I) matching the style of row 6 ---
$ws.Range("A6:I6").Copy() | Out-Null
$ws.Range("A9:I9").PasteSpecial(-4122) | Out-Null

# --- Row 10: new header row (A:I) matching the style of row 7 ---
$ws.Range("A7:I7").Copy() | Out-Null
$ws.Range("A10:I10").PasteSpecial(-4122) | Out-Null
$ws.Range("A10").Value = "Network"
$ws.Range("B10").Value = "Solver"
$ws.Range("C10").Value = "Max epochs"
$ws.Range("D10").Value = "Learning rate"
$ws.Range("E10").Value = "Mini Batch Size"
$ws.Range("F10").Value = "Validation Frequency"
$ws.Range("G10").Value = "Accuracy"
$ws.Range("H10").Value = "Specificity"
$ws.Range("I10").Value = "Sensitivity"

# --- Row 11: new ResNet18 data row (A:I) matching the style of row 5 ---
$ws.Range("A5:I5").Copy() | Out-Null
$ws.Range("A11:I11").PasteSpecial(-4122) | Out-Null
$ws.Range("A11").Value = "ResNet18"
$ws.Range("B11").Value = "Adam"
$ws.Range("C11").Value = 20
$ws.Range("D11").Value = 0.0001
$ws.Range("E11").Value = 64
$ws.Range("F11").Value = 5
$ws.Range("G11").Value = 0.9652
$ws.Range("H11").Value = 0.9766
$ws.Range("I11").Value = 0.9688

# --- Row 12: new "spacer" row (A:I) matching the style of row 6 ---
$ws.Range("A6:I6").Copy() | Out-Null
$ws.Range("A12:I12").PasteSpecial(-4122) | Out-Null

# --- Row 13: new header row (A:I) matching the style of row 7 ---
$ws.Range("A7:I7").Copy() | Out-Null
$ws.Range("A13:I13").PasteSpecial(-4122) | Out-Null
$ws.Range("A13").Value = "Network"
$ws.Range("B13").Value = "Solver"
$ws.Range("C13").Value = "Max epochs"
$ws.Range("D13").Value = "Learning rate"
$ws.Range("E13").Value = "Mini Batch Size"
$ws.Range("F13").Value = "Validation Frequency"
$ws.Range("G13").Value = "Accuracy"
$ws.Range("H13").Value = "Specificity"
$ws.Range("I13").Value = "Sensitivity"

# --- Row 14: new ResNet101 data row (A:I) matching the style of row 5 ---
$ws.Range("A5:I5").Copy() | Out-Null
$ws.Range("A14:I14").PasteSpecial(-4122) | Out-Null
$ws.Range("A14").Value = "ResNet101"
$ws.Range("B14").Value = "Adam"
$ws.Range("C14").Value = 20
$ws.Range("D14").Value = 0.0001
$ws.Range("E14").Value = 64
$ws.Range("F14").Value = 5
$ws.Range("G14").Value = 0.9652
$ws.Range("H14").Value = 0.9682
$ws.Range("I14").Value = 0.9689

$excel.CutCopyMode = $false

# --- Update the selection to match the recorded active cell ---
$ws.Range("J14").Select() | Out-Null
